$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new BOM row (row 23): Programming header / part number / datasheet link / qty
$ws.Range("A23").Value = "Programming header"
$ws.Range("B23").Value = 15910100
$ws.Range("C23").Value = "https://www.digikey.com/en/products/detail/molex/0015910100/614775"
$ws.Range("D23").Value = 1

# Update selection to reflect where the user clicked next (below the new row)
$ws.Range("B24").Select()
